$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.755.28"
$ws.Range("E2").Value = "  -2.90%  "
$ws.Range("D3").Value = "3.172.77"
$ws.Range("E3").Value = "  -1.88%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'599.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("D6").Value = "'151.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.32%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "3.169.27"
$ws.Range("E8").Value = "  -1.92%  "
$ws.Range("E9").Value = "  -3.73%  "
$ws.Range("E10").Value = "  -5.35%  "
$ws.Range("D11").Value = "'5.58"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.13%  "
$ws.Range("E12").Value = "  -6.36%  "
$ws.Range("E13").Value = "  -5.65%  "
$ws.Range("D14").Value = "'36.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.75%  "
$ws.Range("D15").Value = "3.693.17"
$ws.Range("E15").Value = "  -1.90%  "
$ws.Range("D16").Value = "64.743.33"
$ws.Range("E16").Value = "  -3.00%  "
$ws.Range("D17").Value = "3.171.89"
$ws.Range("E17").Value = "  -1.65%  "
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("D19").Value = "'7.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.08%  "
$ws.Range("D20").Value = "'480.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.21%  "
$ws.Range("E21").Value = "  -2.86%  "
$ws.Range("E22").Value = "  -2.89%  "
$ws.Range("D23").Value = "'7.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.05%  "
$ws.Range("D24").Value = "'13.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.59%  "
$ws.Range("D25").Value = "'84.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.19%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "'2.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.87%  "
$ws.Range("D28").Value = "'8.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.47%  "
$ws.Range("E29").Value = "  -5.90%  "
$ws.Range("D30").Value = "'7.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("D31").Value = "'0.120"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.83%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.21%  "
$ws.Range("E33").Value = "  -8.23%  "
$ws.Range("D34").Value = "'26.85"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.88%  "
$ws.Range("E35").Value = "  -6.01%  "
$ws.Range("D36").Value = "'6.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.98%  "
$ws.Range("D37").Value = "'54.67"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.61%  "
$ws.Range("D38").Value = "'3.19"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.57%  "
$ws.Range("E39").Value = "  -5.23%  "
$ws.Range("D40").Value = "'458.54"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.11%  "
$ws.Range("E41").Value = "  -3.09%  "
$ws.Range("D42").Value = "'0.0403"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.86%  "
$ws.Range("D43").Value = "'8.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.96%  "
$ws.Range("E44").Value = "  -1.86%  "
$ws.Range("D45").Value = "2.875.90"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("E46").Value = "  -8.16%  "
$ws.Range("D47").Value = "'27.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.05%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("E49").Value = "  -0.72%  "
$ws.Range("E50").Value = "  -3.51%  "
$ws.Range("D51").Value = "'119.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.25%  "
